function Test-ExactMatch($s1, $s2) {
    if ($s1.Length -ne $s2.Length) { return $false }
    $arr1 = $s1.ToCharArray()
    $arr2 = $s2.ToCharArray()
    for ($i = 0; $i -lt $arr1.Length; $i++) {
        if ([int]$arr1[$i] -ne [int]$arr2[$i]) { return $false }
    }
    return $true
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ([string]::IsNullOrEmpty($val)) {
        continue
    }

    $parts = $val -split ', '

    $hasSystem = $false
    foreach ($p in $parts) {
        if (Test-ExactMatch $p 'System') { $hasSystem = $true }
    }

    if ($hasSystem) {
        $rest = @()
        foreach ($p in $parts) {
            if (-not (Test-ExactMatch $p 'System')) {
                $rest += $p
            }
        }
        $newParts = @('System') + $rest
        $newVal = $newParts -join ', '
        $cell.Value2 = $newVal
    }
}
